# Updates the numeric experiment results to the "alpha_zero" non-convex
# re-run values, per commit "expermits todos no convexos menos el 5to".
#
# All of the values in this workbook (even the numeric-looking ones, e.g.
# "-8.05") are stored as plain text, not as numbers, so every write below
# is forced to text (leading "'" quote prefix) and the style that the
# quote-prefix registers is immediately reset back to Normal so the cell
# keeps the workbook's default (unstyled) formatting.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, [string]$value)
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Worksheet names "Vector_bf" and "Vector_BF" only differ by case, and
# Worksheets.Item() lookups are case-insensitive, so they must be
# addressed by their (1-based) tab position instead of by name.
$wsFuncObj  = $wb.Worksheets.Item(1)   # Funciones_Objetivo        (unchanged)
$wsLider    = $wb.Worksheets.Item(2)   # Restricciones_del_lider   (unchanged)
$wsFollower = $wb.Worksheets.Item(3)   # Restricciones_del_follower
$wsPunto    = $wb.Worksheets.Item(4)   # Punto_modificado
$wsVecbf    = $wb.Worksheets.Item(5)   # Vector_bf
$wsVecBF    = $wb.Worksheets.Item(6)   # Vector_BF
$wsAlpha    = $wb.Worksheets.Item(7)   # Vector_Alpha              (unchanged)

# --- Restricciones_del_follower: rows 2-5 (Expression, Function_Evaluation,
#     Restriction_Set_Type, Lambda_value, Beta_value, Gamma_value) ---
Set-TextValue $wsFollower.Range("A2") "8.05 - y"
Set-TextValue $wsFollower.Range("B2") "-8.05"
Set-TextValue $wsFollower.Range("C2") "J_0_L0_v"
Set-TextValue $wsFollower.Range("D2") "0.13"
Set-TextValue $wsFollower.Range("E2") "8.7"
Set-TextValue $wsFollower.Range("F2") "9.5"

Set-TextValue $wsFollower.Range("A3") "-1.950000000000001 - x + y"
Set-TextValue $wsFollower.Range("B3") "-1.049999999999999"
Set-TextValue $wsFollower.Range("C3") "J_0_L0_v"
Set-TextValue $wsFollower.Range("D3") "0.6"
Set-TextValue $wsFollower.Range("E3") "-0.8"
Set-TextValue $wsFollower.Range("F3") "-0.7000000000000001"

Set-TextValue $wsFollower.Range("A4") "-22.200000000000003 + x + 2y"
Set-TextValue $wsFollower.Range("B4") "10.200000000000001"
Set-TextValue $wsFollower.Range("C4") "J_0_LP_v"
Set-TextValue $wsFollower.Range("D4") "0.62"
Set-TextValue $wsFollower.Range("E4") "9.1"
Set-TextValue $wsFollower.Range("F4") "0"

Set-TextValue $wsFollower.Range("A5") "-17.049999999999997 + 4x - y"
Set-TextValue $wsFollower.Range("B5") "4.349999999999998"
Set-TextValue $wsFollower.Range("C5") "J_Ne_L0_v"
Set-TextValue $wsFollower.Range("D5") "0.32"
Set-TextValue $wsFollower.Range("E5") "-3.8"
Set-TextValue $wsFollower.Range("F5") "-5.2"

# --- Punto_modificado: row 2 (x, y) ---
Set-TextValue $wsPunto.Range("A2") "6.1"
Set-TextValue $wsPunto.Range("B2") "8.05"

# --- Vector_bf: row 2 ---
Set-TextValue $wsVecbf.Range("A2") "-2.3899999999999997"

# --- Vector_BF: rows 2-3 ---
Set-TextValue $wsVecBF.Range("A2") "6.299999999999999"
Set-TextValue $wsVecBF.Range("A3") "-9.5"

# Vector_Alpha is unchanged (stays 0).
